$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(210).Insert()

$ws.Cells.Item(210, 1).Value = 11
$ws.Cells.Item(210, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value = "Bíobío"
$ws.Cells.Item(210, 4).Value = 44889
$ws.Cells.Item(210, 5).Value = 8
$ws.Cells.Item(210, 6).Value = 100112040
$ws.Cells.Item(210, 7).Value = "Cilantro"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 220
$ws.Cells.Item(210, 11).Value = 15000
$ws.Cells.Item(210, 12).Value = 16000
$ws.Cells.Item(210, 13).Value = 15545
$ws.Cells.Item(210, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(210, 15).Value = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value = 432
$ws.Cells.Item(210, 17).Value = 36
$ws.Cells.Item(210, 18).Value = "Hortaliza"
